# Trade #56 closed at 2026-02-17 15:43:25 - unknown UNKNOWN +0.000%
#
# Updates the Summary + Strategy Status roll-up numbers for the
# MarketMaking strategy's new closed trade, and appends the trade's row
# to both the "All Trades" and "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.3    # Current Capital
$summary.Range("B4").Value = 0.3       # Total P&L $
$summary.Range("B5").Value = 0.11      # Total P&L %
$summary.Range("B6").Value = 56        # Total Trades
$summary.Range("B8").Value = 31        # Losing Trades
$summary.Range("B9").Value = 28.57     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.3      # Capital
$status.Range("D4").Value = 56         # Trades
$status.Range("E4").Value = 0.3        # P&L $
$status.Range("F4").Value = 0.3        # P&L %
$status.Range("G4").Value = 28.57      # Win Rate %

# ---------------------------------------------------------------------
# 3. Append new trade row (#56) to "All Trades" and "MarketMaking" logs
# ---------------------------------------------------------------------
$tradeRow = 57

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A$tradeRow").Value = 56
    $ws.Range("A$tradeRow").ClearFormats()

    # The Date column looks like a date to the auto-detect logic, so force
    # it to Text first, then strip the resulting format so the cell keeps
    # the plain (unstyled) look the rest of the sheet uses.
    $ws.Range("B$tradeRow").NumberFormat = "@"
    $ws.Range("B$tradeRow").Value = "2026-02-17"
    $ws.Range("B$tradeRow").ClearFormats()

    $ws.Range("C$tradeRow").Value = "15:43:18"
    $ws.Range("C$tradeRow").ClearFormats()

    $ws.Range("D$tradeRow").Value = "MarketMaking"
    $ws.Range("D$tradeRow").ClearFormats()

    $ws.Range("E$tradeRow").Value = "UP"
    $ws.Range("E$tradeRow").ClearFormats()

    $ws.Range("F$tradeRow").Value = 0.27
    $ws.Range("F$tradeRow").ClearFormats()

    $ws.Range("G$tradeRow").Value = 0.2
    $ws.Range("G$tradeRow").ClearFormats()

    $ws.Range("H$tradeRow").Value = "CLOSED"
    $ws.Range("H$tradeRow").ClearFormats()

    $ws.Range("I$tradeRow").Value = -25.9259
    $ws.Range("I$tradeRow").ClearFormats()

    $ws.Range("J$tradeRow").Value = -0.07000000000000001
    $ws.Range("J$tradeRow").ClearFormats()

    $ws.Range("K$tradeRow").Value = 100.3
    $ws.Range("K$tradeRow").ClearFormats()

    $ws.Range("L$tradeRow").Value = 0
    $ws.Range("L$tradeRow").ClearFormats()

    $ws.Range("M$tradeRow").Value = 0
    $ws.Range("M$tradeRow").ClearFormats()

    $ws.Range("N$tradeRow").Value = 0.6
    $ws.Range("N$tradeRow").ClearFormats()

    $ws.Range("O$tradeRow").Value = "Normal spread capture: 19600 bps"
    $ws.Range("O$tradeRow").ClearFormats()

    $ws.Range("P$tradeRow").Value = "early_exit"
    $ws.Range("P$tradeRow").ClearFormats()

    $ws.Range("Q$tradeRow").Value = 0.13
    $ws.Range("Q$tradeRow").ClearFormats()
}
